# Add a "Normal weight" category column to the LDL study summary sheet.
# This inserts a new column S ("Normal weight", Yes/No per study) before the
# existing "Link" column (which shifts from S to T), updates the three
# existing hyperlinks so they point at their new T-column locations, and
# moves the sheet selection / dimension to reflect the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: detach the hyperlink relationships from their current (soon to
# be stale) S-column cells so re-inserting the column doesn't leave behind
# duplicate / orphaned hyperlink entries. Cell text/value/style are untouched.
$ws.Range("S12").Hyperlinks.Delete()
$ws.Range("S32").Hyperlinks.Delete()
$ws.Range("S33").Hyperlinks.Delete()

# --- Step 2: insert a new blank column at S (column 19), shifting the old
# S column ("Link") - and everything after the insertion point - one column
# to the right, into T.
$ws.Columns.Item(19).Insert()

# --- Step 3: populate the new "Normal weight" column (header + per-study
# Yes/No values) for every row that has data for it.
$ws.Range("S1").Value = "Normal weight"
$ws.Range("S2").Value = "No"
$ws.Range("S3").Value = "No"
$ws.Range("S4").Value = "No"
$ws.Range("S7").Value = "No"
$ws.Range("S10").Value = "No"
$ws.Range("S11").Value = "No"
$ws.Range("S12").Value = "No"
$ws.Range("S13").Value = "Yes"
$ws.Range("S15").Value = "No"
$ws.Range("S20").Value = "Yes"
$ws.Range("S21").Value = "Yes"
$ws.Range("S24").Value = "Yes"
$ws.Range("S30").Value = "No"
$ws.Range("S31").Value = "No"
$ws.Range("S32").Value = "No"
$ws.Range("S33").Value = "No"

# --- Step 4: re-create the three hyperlinks at their new T-column homes,
# then restore the standard Excel "Hyperlink" cell style (blue/underline)
# that was already carried over from the column shift, so the visible
# formatting matches the original.
$ws.Hyperlinks.Add($ws.Range("T12"), "https://pubmed.ncbi.nlm.nih.gov/23155696/")
$ws.Range("T12").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("T32"), "https://academic.oup.com/jcem/article/88/4/1617/2845298?login=false")
$ws.Range("T32").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("T33"), "https://www.nejm.org/doi/10.1056/NEJMoa022207?url_ver=Z39.88-2003&rfr_id=ori:rid:crossref.org&rfr_dat=cr_pub%20%200www.ncbi.nlm.nih.gov")
$ws.Range("T33").Style = "Hyperlink"

# --- Step 5: update the sheet selection to match the saved workbook state.
$ws.Range("T39").Select()
